$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append, corresponding to the new rows 372..379 in the diff.
$rows = @(
    @("Create a country",     "PASSED", "chrome"),
    @("Create a country 2",   "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome")
)

$startRow = 372
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
